$d = $word.ActiveDocument

# The paragraph currently reads "Versi" + "on" + " 2" + "." (separate runs,
# spelling across two runs as "Version" with a proofed-word wrapper) and
# needs to become "Version" + " 1." (i.e. "Version 1." overall), dropping
# the now-superfluous trailing "." run.

# Step 1: the word is split across two runs ("Versi" / "on"). Collapse them
# into a single run by deleting the text of the "on" run; the remaining
# "Versi" run will be renamed to "Version" next.
$onRange = $d.Content
$onRange.Find.Execute("on") | Out-Null
$onRange.Text = ""

# Step 2: rename "Versi" -> "Version".
$d.Content.Find.Execute("Versi", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Version", 2)

# Step 3: the " 2" run becomes " 1." (version number bump + moved period).
$d.Content.Find.Execute("2", $false, $false, $false, $false, $false, `
    $true, 1, $false, "1.", 2)

# Step 4: the text now reads "Version 1.." -- the final character is the
# old standalone "." run, which is redundant now that the period moved
# into the " 1." run. Locate the ".." and drop the second (trailing) dot.
$dotdot = $d.Content
$dotdot.Find.Execute("..") | Out-Null
$trailingDot = $d.Range($dotdot.Start + 1, $dotdot.End)
$trailingDot.Text = ""
